$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: cell reference -> new text value
$updates = @(
    @{ Cell = 'D2'; Value = '58.481.70' }
    @{ Cell = 'E2'; Value = '  -0.01%  ' }
    @{ Cell = 'D3'; Value = '2.634.40' }
    @{ Cell = 'E3'; Value = '  +3.27%  ' }
    @{ Cell = 'E4'; Value = '  +0.16%  ' }
    @{ Cell = 'D5'; Value = '510.40' }
    @{ Cell = 'E5'; Value = '  +1.01%  ' }
    @{ Cell = 'D6'; Value = '142.20' }
    @{ Cell = 'E6'; Value = '  -1.08%  ' }
    @{ Cell = 'E7'; Value = '  -0.28%  ' }
    @{ Cell = 'D8'; Value = '0.562' }
    @{ Cell = 'E8'; Value = '  +1.93%  ' }
    @{ Cell = 'D9'; Value = '2.669.51' }
    @{ Cell = 'E9'; Value = '  +4.48%  ' }
    @{ Cell = 'D10'; Value = '6.26' }
    @{ Cell = 'E10'; Value = '  +1.28%  ' }
    @{ Cell = 'E11'; Value = '  +3.32%  ' }
    @{ Cell = 'D12'; Value = '0.333' }
    @{ Cell = 'E12'; Value = '  +0.69%  ' }
    @{ Cell = 'E13'; Value = '  -1.31%  ' }
    @{ Cell = 'D14'; Value = '3.106.43' }
    @{ Cell = 'E14'; Value = '  +3.59%  ' }
    @{ Cell = 'D15'; Value = '58.533.89' }
    @{ Cell = 'E15'; Value = '  +0.08%  ' }
    @{ Cell = 'D16'; Value = '20.74' }
    @{ Cell = 'E16'; Value = '  +1.27%  ' }
    @{ Cell = 'D17'; Value = '0.0000136' }
    @{ Cell = 'E17'; Value = '  +2.00%  ' }
    @{ Cell = 'D18'; Value = '2.665.68' }
    @{ Cell = 'E18'; Value = '  +4.28%  ' }
    @{ Cell = 'D19'; Value = '4.51' }
    @{ Cell = 'E19'; Value = '  +0.06%  ' }
    @{ Cell = 'D20'; Value = '340.62' }
    @{ Cell = 'E20'; Value = '  +2.24%  ' }
    @{ Cell = 'D21'; Value = '10.34' }
    @{ Cell = 'E21'; Value = '  +3.02%  ' }
    @{ Cell = 'D22'; Value = '6.08' }
    @{ Cell = 'E22'; Value = '  +2.52%  ' }
    @{ Cell = 'D23'; Value = '0.999' }
    @{ Cell = 'E23'; Value = '  +0.21%  ' }
    @{ Cell = 'D24'; Value = '60.75' }
    @{ Cell = 'E24'; Value = '  +2.33%  ' }
    @{ Cell = 'D25'; Value = '0.417' }
    @{ Cell = 'E25'; Value = '  +2.59%  ' }
    @{ Cell = 'D26'; Value = '2.769.06' }
    @{ Cell = 'E26'; Value = '  +3.60%  ' }
    @{ Cell = 'D27'; Value = '0.995' }
    @{ Cell = 'E27'; Value = '  -0.62%  ' }
    @{ Cell = 'D28'; Value = '0.159' }
    @{ Cell = 'E28'; Value = '  +2.56%  ' }
    @{ Cell = 'D29'; Value = '0.0₃0799' }
    @{ Cell = 'E29'; Value = '  +3.36%  ' }
    @{ Cell = 'D30'; Value = '7.14' }
    @{ Cell = 'E30'; Value = '  +4.31%  ' }
    @{ Cell = 'D31'; Value = '0.997' }
    @{ Cell = 'E31'; Value = '  -0.26%  ' }
    @{ Cell = 'D32'; Value = '6.32' }
    @{ Cell = 'E32'; Value = '  +8.15%  ' }
    @{ Cell = 'D33'; Value = '18.81' }
    @{ Cell = 'E33'; Value = '  +1.51%  ' }
    @{ Cell = 'D34'; Value = '1.56' }
    @{ Cell = 'E34'; Value = '  +1.74%  ' }
    @{ Cell = 'D35'; Value = '149.31' }
    @{ Cell = 'E35'; Value = '  +0.35%  ' }
    @{ Cell = 'D36'; Value = '1.00' }
    @{ Cell = 'E36'; Value = '  +10.19%  ' }
    @{ Cell = 'D37'; Value = '3.97' }
    @{ Cell = 'E37'; Value = '  +2.39%  ' }
    @{ Cell = 'D38'; Value = '1.13' }
    @{ Cell = 'E38'; Value = '  +2.91%  ' }
    @{ Cell = 'B39'; Value = 'OKB' }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb' }
    @{ Cell = 'D39'; Value = '36.49' }
    @{ Cell = 'E39'; Value = '  +1.73%  ' }
    @{ Cell = 'B40'; Value = 'Fetch.AI' }
    @{ Cell = 'C40'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' }
    @{ Cell = 'D40'; Value = '0.840' }
    @{ Cell = 'E40'; Value = '  +3.01%  ' }
    @{ Cell = 'D41'; Value = '3.64' }
    @{ Cell = 'E41'; Value = '  +3.60%  ' }
    @{ Cell = 'D42'; Value = '1.39' }
    @{ Cell = 'E42'; Value = '  +1.22%  ' }
    @{ Cell = 'B43'; Value = 'Mantle' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt' }
    @{ Cell = 'D43'; Value = '0.615' }
    @{ Cell = 'E43'; Value = '  +1.48%  ' }
    @{ Cell = 'B44'; Value = 'Bittensor' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' }
    @{ Cell = 'D44'; Value = '277.83' }
    @{ Cell = 'E44'; Value = '  -1.82%  ' }
    @{ Cell = 'D45'; Value = '0.995' }
    @{ Cell = 'E45'; Value = '  -0.39%  ' }
    @{ Cell = 'D46'; Value = '0.0973' }
    @{ Cell = 'E46'; Value = '  -0.72%  ' }
    @{ Cell = 'D47'; Value = '19.50' }
    @{ Cell = 'E47'; Value = '  +4.99%  ' }
    @{ Cell = 'D48'; Value = '0.0528' }
    @{ Cell = 'E48'; Value = '  -0.58%  ' }
    @{ Cell = 'B49'; Value = 'WhiteBITCoin' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt' }
    @{ Cell = 'D49'; Value = '10.26' }
    @{ Cell = 'E49'; Value = '  -0.65%  ' }
    @{ Cell = 'B50'; Value = 'VeChain' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D50'; Value = '0.0228' }
    @{ Cell = 'E50'; Value = '  +1.23%  ' }
    @{ Cell = 'D51'; Value = '4.68' }
    @{ Cell = 'E51'; Value = '  +4.15%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Cell.Substring(0,1) -eq "D") {
        # Force column D (Price) to remain plain text so values like
        # "1.00", "0.0000136" or "58.481.70" are not reinterpreted as numbers/dates.
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
